$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "LogFilePath" before the existing "strEmailAccount" row (row 4).
$ws.Rows.Item(4).Insert()

# Insert a new row for "WorkpackageName" after the "strEmailAccount" row (now row 5).
$ws.Rows.Item(6).Insert()

# Fill in the new LogFilePath row (row 4).
$ws.Range("A4").Value = "LogFilePath"
$ws.Range("B4").Value = "C:\Users\{0}\Desktop\Demo Robot Log_{1}.xlsx"
$ws.Range("C4").Value = "log file path"

# Fill in the new WorkpackageName row (row 6).
$ws.Range("A6").Value = "WorkpackageName"
$ws.Range("B6").Value = "Demo Robot"
$ws.Range("C6").Value = "workpackage name"

# Resize the table to include the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C8"))

# Update the saved selection to match the final workbook state.
[void]$ws.Range("B19").Select()
